# Modifying esmith10laptop's Dropbox folder.
# The esmith10laptop column (I) used to point at a SkyDrive-based
# LivemRNAData path; it now uses a plain Dropbox folder instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# DropboxFolder (row 5) for esmith10laptop (column I).
$ws.Range("I5").Value = "C:\E\Dropbox\LivemRNAData"

# Reflect the edited cell becoming the active selection / split-pane
# scroll position, as captured by the author's saved view state.
$excel.ActiveWindow.SplitColumn = 7
$ws.Range("I5").Select()
